$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - this also resets sharedStrings so we fully control their order.
$ws.Cells.Clear()

# --- Seed shared strings in the exact order the target file expects ---
# (m3/s) first (was already first among the surviving strings), then the
# 8 plant names in row order, then the brand-new header labels.
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("C2").Value = "Perlen 1 (WTA HF)"
$ws.Range("C3").Value = "Perlen 2 (WTA PF)"
$ws.Range("C4").Value = "Thorenberg"
$ws.Range("C5").Value = "Rathausen"
$ws.Range("C6").Value = "Wolhusen (Geistlich)"
$ws.Range("C7").Value = "Emmenweid"
$ws.Range("C8").Value = "Mühlenplatz"
$ws.Range("C9").Value = "Stollen"
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# --- Header row styling: F1:K1 need font 9pt Arial, no explicit number
# format applied (distinct from the body "s=1" style which does mark
# applyNumberFormat). Building a throwaway named style and deleting it
# afterwards yields that exact xf shape without leaving cellStyle debris.
$wb.Styles.Add("TmpHeaderStyle")
$hs = $wb.Styles.Item("TmpHeaderStyle")
$hs.Font.Name = "Arial"
$hs.Font.Size = 9
$ws.Range("F1:K1").Style = "TmpHeaderStyle"
$wb.Styles.Item("TmpHeaderStyle").Delete()

# --- Data rows (A2:K9) ---
$data = @(
    @(1, 304600, "Perlen 1 (WTA HF)", 1873, 1981, 45, 1, 0.96, 4, 4, 8),
    @(2, 304700, "Perlen 2 (WTA PF)", 1875, 2000, 45, 1.1499999999999999, 1.0900000000000001, 3.9, 3.9, 7.8),
    @(3, 304300, "Thorenberg", 1886, 2000, 7, 0.79, 0.75, 2.2000000000000002, 2.6, 4.8),
    @(4, 304500, "Rathausen", 1896, 1980, 45, 2.1, 2, 7.88, 8.02, 15.9),
    @(5, 304200, "Wolhusen (Geistlich)", 1906, 2003, 4, 0.33, 0.3, 0.53, 0.63, 1.1599999999999999),
    @(6, 304400, "Emmenweid", 1931, 2003, 12, 1.08, 1.08, 1.59, 2.96, 4.55),
    @(7, 304150, "Mühlenplatz", 1998, $null, 58, 0.93, 0.83, 2, 1, 3),
    @(8, 304250, "Stollen", 1999, $null, 0.35, 0.71, 0.64, 1, 2, 3)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 1).Font.Name = "Arial"
    $ws.Cells.Item($row, 1).Font.Size = 9
    $ws.Cells.Item($row, 1).NumberFormat = "0"

    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 2).Font.Name = "Arial"
    $ws.Cells.Item($row, 2).Font.Size = 9
    $ws.Cells.Item($row, 2).NumberFormat = "0"

    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 3).Font.Name = "Arial"
    $ws.Cells.Item($row, 3).Font.Size = 9

    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 4).Font.Name = "Arial"
    $ws.Cells.Item($row, 4).Font.Size = 9
    $ws.Cells.Item($row, 4).NumberFormat = "0"

    if ($r[4] -ne $null) {
        $ws.Cells.Item($row, 5).Value = $r[4]
        $ws.Cells.Item($row, 5).Font.Name = "Arial"
        $ws.Cells.Item($row, 5).Font.Size = 9
        $ws.Cells.Item($row, 5).NumberFormat = "0"
    }

    for ($col = 6; $col -le 11; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $r[$col - 1]
        $cell.Font.Name = "Arial"
        $cell.Font.Size = 9
        $cell.NumberFormat = "0.00"
    }

    $row = $row + 1
}

# --- View bits: selection + active cell ---
$ws.Range("A2:K2").Select()

Write-Host "Applied LU 2006 canton restructuring"
